# Add 2022-Q4 data.
#
# Before:  总计, 2022-Q3, 2021-Q3
# After:   总计, 2022-Q4, 2022-Q3, 2021-Q3
#
# 1. A brand new worksheet "2022-Q4" is inserted right after "总计" and
#    before "2022-Q3". It reuses the same layout/headers as the other
#    quarterly detail sheets, but with refreshed fund figures.
# 2. The "总计" summary sheet gets a new row for "2022-Q4" inserted right
#    after the existing (now outdated) 2022-Q3-labelled row, which itself
#    becomes the "2022-Q3" row, and the old "2021-Q3" row shifts down.

$wb = $excel.ActiveWorkbook

$totalSheet = $wb.Worksheets.Item(1)
$q3Sheet = $wb.Worksheets.Item(2)

# --- 1. Insert the new "2022-Q4" worksheet right before "2022-Q3" ---
$q4Sheet = $wb.Worksheets.Add($q3Sheet)
$q4Sheet.Name = "2022-Q4"

# Header row, formatted like the other quarterly sheets (bold, centered,
# boxed).
$q4Sheet.Range("B1").Value = "基金代码"
$q4Sheet.Range("C1").Value = "基金名称"
$q4Sheet.Range("D1").Value = "基金规模"
$q4Sheet.Range("E1").Value = "股票总仓位"
$q4Sheet.Range("F1").Value = "仓位占比"
$q4Sheet.Range("G1").Value = "持有市值(亿元)"
$q4Sheet.Range("H1").Value = "仓位排名"

$headerRange = $q4Sheet.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Index column (A) uses the same bold/boxed/centered look as the header.
$indexRange = $q4Sheet.Range("A2:A3")
$indexRange.Font.Bold = $true
$indexRange.HorizontalAlignment = -4108
$indexRange.VerticalAlignment = -4160
$indexRange.Borders.LineStyle = 1

# Row 2 - 泰康景泰回报混合A (values refreshed for 2022-Q4)
$q4Sheet.Range("A2").Value = 0
$q4Sheet.Range("B2").NumberFormat = "@"
$q4Sheet.Range("B2").Value = "005014"
$q4Sheet.Range("C2").Value = "泰康景泰回报混合A"
$q4Sheet.Range("D2").NumberFormat = "@"
$q4Sheet.Range("D2").Value = "8.86"
$q4Sheet.Range("E2").NumberFormat = "@"
$q4Sheet.Range("E2").Value = "32.91"
$q4Sheet.Range("F2").NumberFormat = "@"
$q4Sheet.Range("F2").Value = "1.23"
$q4Sheet.Range("G2").NumberFormat = "@"
$q4Sheet.Range("G2").Value = "0.1090"
$q4Sheet.Range("H2").Value = 8

# Row 3 - 泰康景泰回报混合C (values refreshed for 2022-Q4)
$q4Sheet.Range("A3").Value = 1
$q4Sheet.Range("B3").NumberFormat = "@"
$q4Sheet.Range("B3").Value = "005015"
$q4Sheet.Range("C3").Value = "泰康景泰回报混合C"
$q4Sheet.Range("D3").NumberFormat = "@"
$q4Sheet.Range("D3").Value = "0.37"
$q4Sheet.Range("E3").NumberFormat = "@"
$q4Sheet.Range("E3").Value = "32.91"
$q4Sheet.Range("F3").NumberFormat = "@"
$q4Sheet.Range("F3").Value = "1.23"
$q4Sheet.Range("G3").NumberFormat = "@"
$q4Sheet.Range("G3").Value = "0.0046"
$q4Sheet.Range("H3").Value = 8

# --- 2. Update the "总计" summary sheet ---
# Shift the old "2021-Q3" row (row 3) down to row 4, keep its values.
$totalSheet.Range("A4").Value = 2
$totalSheet.Range("A4").Font.Bold = $true
$totalSheet.Range("A4").HorizontalAlignment = -4108
$totalSheet.Range("A4").VerticalAlignment = -4160
$totalSheet.Range("A4").Borders.LineStyle = 1
$totalSheet.Range("B4").Value = "2021-Q3"
$totalSheet.Range("C4").Value = 1
$totalSheet.Range("D4").Value = 0.19

# Row 3 becomes the "2022-Q3" row, carrying the same totals that used to
# sit in row 2 (count of holdings / market value were unchanged between
# Q3 and Q4).
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("B3").Value = "2022-Q3"
$totalSheet.Range("C3").Value = 2
$totalSheet.Range("D3").Value = 0.11

# Row 2 now reports the new "2022-Q4" totals.
$totalSheet.Range("B2").Value = "2022-Q4"

# Restore "总计" as the active sheet/tab (inserting a worksheet makes the
# newly inserted one active by default).
$totalSheet.Activate()
